$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.785.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.628.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.854.61"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.624.39"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.551"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.792.43"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.60"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.28"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.996"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.18"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.49"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.902"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.140.33"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.49"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.995"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.34"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.799"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.765.44"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.37"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.46"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.73%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.31%  "
